{"js": "// Add two paragraphs of \"Critical Reflection\" text.\n// Both target paragraphs already exist in the document (empty paragraphs,\n// one right after the \"Encountered difficulties\" discussion, the other\n// right after the \"Improvement/Reflection\" heading) \u2014 we only need to\n// fill them in with the authored text, using the same font as the rest\n// of the document body (\"Franklin Gothic Book\").\n\nconst FONT = \"Franklin Gothic Book\";\n\nconst lcaText =\n  \"The lowest common ancestor was difficult to solve, it consists of several parts. \" +\n  \"Storing the parent of each node, if the children are the nodes that are being searched for return the parent. \" +\n  \"If the nodes searched for are in separate paths from the root node return the root node. \" +\n  \"The challenging problem is in returning the parent of two nodes who don\\u2019t share the same parent, \" +\n  \"this requires the algorithm to back track to find the lowest common ancestor one way of doing this would be to store the paths to each node and compare them to find the LCA of them, \" +\n  \"another would be to backtrack till a common ancestor is found. \" +\n  \"The choice was the former, this meant comparing the two paths, however, this was not successfully implemented \" +\n  \"as any trees which span more than 4 depths were difficult to assess and to return the LCA.\";\n\nconst reflectionText =\n  \"The task overall worked well, however, would have correctly configured a function to \" +\n  \"correctly return the LCA of two which are divergent of different parents. This works on a \" +\n  \"small-scale tree but a better algorithm to correctly return the LCA would be better. \" +\n  \"Backtracking and identifying would have done this. This was difficult to visualise and then to code. \" +\n  \"All in all the task was completed with some error catching however would like to add some filtration \" +\n  \"for user input of entries in the tree to search for.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"uniqueLocalId\");\n}\nawait context.sync();\n\nlet lcaParagraph = null;\nlet reflectionParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const id = paragraphs.items[i].uniqueLocalId;\n  if (id === \"2CD052B4\") {\n    lcaParagraph = paragraphs.items[i];\n  } else if (id === \"18540E2B\") {\n    reflectionParagraph = paragraphs.items[i];\n  }\n}\n\nif (!lcaParagraph || !reflectionParagraph) {\n  throw new Error(\"Could not locate target paragraphs for the critical reflection text.\");\n}\n\nconst lcaRange = lcaParagraph.insertText(lcaText, Word.InsertLocation.end);\nlcaRange.font.name = FONT;\n\nconst reflectionRange = reflectionParagraph.insertText(reflectionText, Word.InsertLocation.end);\nreflectionRange.font.name = FONT;\n\nawait context.sync();\n", "ps1": "# Add two paragraphs of \"Critical Reflection\" text.\n# Both target paragraphs already exist in the document (empty paragraphs,\n# one right after the \"Encountered difficulties\" discussion, the other\n# right after the \"Improvement/Reflection\" heading) -- we only need to\n# fill them in with the authored text, using the same font as the rest\n# of the document body (\"Franklin Gothic Book\").\n\n$d = $word.ActiveDocument\n\n$lcaText = \"The lowest common ancestor was difficult to solve, it consists of several parts. \" + `\n    \"Storing the parent of each node, if the children are the nodes that are being searched for return the parent. \" + `\n    \"If the nodes searched for are in separate paths from the root node return the root node. \" + `\n    \"The challenging problem is in returning the parent of two nodes who don\" + [char]0x2019 + \"t share the same parent, \" + `\n    \"this requires the algorithm to back track to find the lowest common ancestor one way of doing this would be to store the paths to each node and compare them to find the LCA of them, \" + `\n    \"another would be to backtrack till a common ancestor is found. \" + `\n    \"The choice was the former, this meant comparing the two paths, however, this was not successfully implemented \" + `\n    \"as any trees which span more than 4 depths were difficult to assess and to return the LCA.\"\n\n$reflectionText = \"The task overall worked well, however, would have correctly configured a function to \" + `\n    \"correctly return the LCA of two which are divergent of different parents. This works on a \" + `\n    \"small-scale tree but a better algorithm to correctly return the LCA would be better. \" + `\n    \"Backtracking and identifying would have done this. This was difficult to visualise and then to code. \" + `\n    \"All in all the task was completed with some error catching however would like to add some filtration \" + `\n    \"for user input of entries in the tree to search for.\"\n\nfunction Get-EmptyParagraphAfter($searchText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($searchText)\n    if (-not $found) {\n        throw \"Could not find anchor text: $searchText\"\n    }\n    # Expand to the whole paragraph that contains the found text, then\n    # collapse to its end -- that lands right at the start of the next\n    # (empty) paragraph.\n    $rng.Expand(4) | Out-Null\n    $rng.Collapse(0) | Out-Null\n    return $rng\n}\n\n$target1 = Get-EmptyParagraphAfter(\"To find the LCA was the most challenging\")\n$target1.InsertAfter($lcaText)\n$target1.Font.Name = \"Franklin Gothic Book\"\n\n$target2 = Get-EmptyParagraphAfter(\"Improvement/Reflection\")\n$target2.InsertAfter($reflectionText)\n$target2.Font.Name = \"Franklin Gothic Book\"\n"}
